$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '247.24'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.36'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.089'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05620'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8135'
$ws.Range("E7").Value = '6MXTokenMX'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8489'
$ws.Range("E8").Value = '7FTXTokenFTT'
$ws.Range("B9").Value = 'BitrueCoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.02821'
$ws.Range("E9").Value = '8BitrueCoinBTR'
$ws.Range("B10").Value = 'BitMartToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09400'
$ws.Range("E10").Value = '9BitMartTokenBMX'
$ws.Range("B11").Value = 'BitForexToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.001515'
$ws.Range("E11").Value = '10BitForexTokenBF'
$ws.Range("B12").Value = 'One'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0005967'
$ws.Range("E12").Value = '11OneONE'
$ws.Range("B13").Value = 'TigerCash'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.006162'
$ws.Range("E13").Value = '12TigerCashTCH'
$ws.Range("B14").Value = 'LEO'
$ws.Range("C14").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.585'
$ws.Range("E14").Value = '13LEOLEO'
$ws.Range("B15").Value = 'GateToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.051'
$ws.Range("E15").Value = '14GateTokenGT'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.1345'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06991'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03169'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1320'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.744'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04661'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1374'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001249'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009595'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001937'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03678'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1361'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002659'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003446'
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008572'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005291'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1199'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002069'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
